$wb = $excel.ActiveWorkbook

# Sheet "展览": update "想去人数" (want-to-go count) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 299
$wsExhibit.Range("F4").Value = 1298

# Sheet "全部类型": update corresponding "想去人数" values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 299
$wsAll.Range("F5").Value = 1298

$wb.Save()
